# "Model 1" (row 4) and "Model 2" (row 5) job items are now finished:
# mark their STATUS as "Done" and fill in the TANGGAL SELESAI (completion date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Completion date: 27 Oct 2021 (matches the other date cells, serial 44496)
$completionDate = 44496

# Row 4 ("Model 1"): set TANGGAL SELESAI (F4) and STATUS (H4)
$ws.Range("F4").NumberFormat = $ws.Range("E4").NumberFormat
$ws.Range("F4").Value = $completionDate
$ws.Range("H4").Value = "Done"

# Row 5 ("Model 2"): set TANGGAL SELESAI (F5) and STATUS (H5)
$ws.Range("F5").NumberFormat = $ws.Range("E4").NumberFormat
$ws.Range("F5").Value = $completionDate
$ws.Range("H5").Value = "Done"

# Active selection ends up on F9 (as last left by the user)
$ws.Range("F9").Select()
